$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at position 14. Excel copies the formatting of the row
#    above (row 13) into the new row, and shifts old rows 14-23 down to 15-24.
# ---------------------------------------------------------------------------
$ws.Rows.Item(14).Insert()

# ---------------------------------------------------------------------------
# 1b. Row 10 ("Objetivos:") gets the new Portuguese objectives paragraph
#     (the old B10/C10 text, "5840535 - Messias Borges Silva", moves down to
#     become the "Docentes responsáveis:" value at row 13).
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = "Fornecer conhecimentos básicos e aplicações das técnicas usuais de manufatura enxuta e estratégia Seis Sigma, apresentando conceitos, princípios e ferramentas utilizados nesta abordagem. Será dado ênfase na mentalidade enxuta, no mapeamento de fluxo de valor, projeto de sistemas de controle de produção enxuta, projeto de células de manufatura, desenvolvimento de trabalhadores multifuncionais e definição de requisitos para a implementação de linhas de produção enxuta."
$ws.Range("C10").Value = "Fornecer conhecimentos básicos e aplicações das técnicas usuais de manufatura enxuta e estratégia Seis Sigma, apresentando conceitos, princípios e ferramentas utilizados nesta abordagem. Será dado ênfase na mentalidade enxuta, no mapeamento de fluxo de valor, projeto de sistemas de controle de produção enxuta, projeto de células de manufatura, desenvolvimento de trabalhadores multifuncionais e definição de requisitos para a implementação de linhas de produção enxuta."

# ---------------------------------------------------------------------------
# 2. Row 13 used to hold "Programa resumido: / Semestral". After the insert it
#    still sits at row 13, but the target layout needs it to hold the
#    "Docentes responsáveis:" value (5840535 - Messias Borges Silva) with no
#    label in column A and default row height.
# ---------------------------------------------------------------------------
$ws.Range("A13").ClearContents()
$ws.Range("B13").Value = "5840535 - Messias Borges Silva"
$ws.Range("C13").Value = "5840535 - Messias Borges Silva"
$ws.Rows.Item(13).RowHeight = 15

# ---------------------------------------------------------------------------
# 3. The freshly inserted row 14 becomes "Programa resumido:" with the new
#    Portuguese summary paragraph.
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B14").Value = "O Pensamento Enxuto; Categorias de Desperdícios; Entendendo o Fluxo de Valor; Mapeamento de Fluxo de Valor; Takt Time; Sistema de Produção Puxada; Criação de fluxo contínuo de produção; Células de Manufatura; Elementos de Controle de Produção; Ferramentas para produção Enxuta. A Estratégia Six Sigma"
$ws.Range("C14").Value = "O Pensamento Enxuto; Categorias de Desperdícios; Entendendo o Fluxo de Valor; Mapeamento de Fluxo de Valor; Takt Time; Sistema de Produção Puxada; Criação de fluxo contínuo de produção; Células de Manufatura; Elementos de Controle de Produção; Ferramentas para produção Enxuta. A Estratégia Six Sigma"
$ws.Rows.Item(14).RowHeight = 60

# ---------------------------------------------------------------------------
# 4. Rows 15-24 (old rows 14-23, shifted down by the insert) already carry the
#    correct column-A labels and row heights; only the B/C content needs to
#    be refreshed where the diff changed it.
# ---------------------------------------------------------------------------

# Row 16 ("Programa:") gets the expanded Portuguese syllabus text.
$ws.Range("B16").Value = "A Abordagem da Produção Enxuta; História; O Pensamento Enxuto; Categorias de Desperdícios; Metodologia DMAIC aplicada ao Lean; Cultura organizacional para o Lean; Entendendo o Fluxo de Valor; Mapeamento de Fluxo de Valor; Fluxo de Material e de Informação; Características do Fluxo de Valor Enxuto; Conceito de Takt Time; Análise de Capacidades; Sistema de Produção Puxada; Criação de fluxo contínuo de produção; Células de Manufatura; Critérios para o Projeto de Células; Diagrama de Espaguete; Preparação da Força de Trabalho; Elementos de Controle de Produção; Kanban; Heijunka Box; Dimensionamento de Kanbans; Aspectos de Gestão da Implantação; Gestão Visual. A Estratégia Seis Sigma e o Método DMAIC."
$ws.Range("C16").Value = "A Abordagem da Produção Enxuta; História; O Pensamento Enxuto; Categorias de Desperdícios; Metodologia DMAIC aplicada ao Lean; Cultura organizacional para o Lean; Entendendo o Fluxo de Valor; Mapeamento de Fluxo de Valor; Fluxo de Material e de Informação; Características do Fluxo de Valor Enxuto; Conceito de Takt Time; Análise de Capacidades; Sistema de Produção Puxada; Criação de fluxo contínuo de produção; Células de Manufatura; Critérios para o Projeto de Células; Diagrama de Espaguete; Preparação da Força de Trabalho; Elementos de Controle de Produção; Kanban; Heijunka Box; Dimensionamento de Kanbans; Aspectos de Gestão da Implantação; Gestão Visual. A Estratégia Seis Sigma e o Método DMAIC."

# Row 19 ("Método:") now holds the "Aulas Expositivas..." text (previously at
# the "Critério:" row).
$ws.Range("B19").Value = "Aulas Expositivas; trabalhos e seminários."
$ws.Range("C19").Value = "Aulas Expositivas; trabalhos e seminários."

# Row 20 ("Critério:") now holds the "MF = ..." grading formula (previously
# at the "Norma de recuperação:" row).
$ws.Range("B20").Value = "MF = (0,30*P1 + 0,30*P2 + 0,40*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários."
$ws.Range("C20").Value = "MF = (0,30*P1 + 0,30*P2 + 0,40*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários."

# Row 21 ("Norma de recuperação:") now holds the "NF = ..." makeup-exam
# formula (previously at the "Bibliografia:" row).
$ws.Range("B21").Value = "NF = (MF + PR)/2, onde PR é uma prova de recuperação"
$ws.Range("C21").Value = "NF = (MF + PR)/2, onde PR é uma prova de recuperação"

# Row 22 ("Bibliografia:") now holds the full bibliography text.
$ws.Range("B22").Value = "WOMACK, J.P.; JONES, D.T. Lean Thinking,  Free Press, 2010WOMACK, J.P.; JONES, D.T., Lean Solutions, 2009DENIS, P, Produçao Lean Simplificada,  Bookman, 2008 HARRY, M. , LINSENMANND.R., The Six Sigma Fieldbook, Doubleday, New York, 2006KUME, H. (Tradução Miyake, D.I.). 1993. Métodos Estatísticos para Melhoria da Qualidade. São Paulo. Editora Gente, 245 p.HARRY, M. , LINSENMANND.R., The Six Sigma Fieldbook, Doubleday, New York, 2006ISHIKAWA, K. Guide to quality control. Tokyo: Kraus Asian Productivity Organization, 1982.JURAN, J.M.; GRYNA, F.M. Juran controle da qualidade: métodos especiais de apoio à qualidade. São Paulo: Makron Books, 1993LIKER, K. & MEIER D. O Modelo Toyota, Manual de Aplicação. Porto Alegre: Bookman, 2007MARANHÃO, M. ISO Série 9000-Guia de Implementação, Qualitymark, Rio de Janeiro, 2001MONTGOMERY, D.C. 1991. Introduction to Statistical Quality Control. New York. John Wiley & Sons Inc."
$ws.Range("C22").Value = "WOMACK, J.P.; JONES, D.T. Lean Thinking,  Free Press, 2010WOMACK, J.P.; JONES, D.T., Lean Solutions, 2009DENIS, P, Produçao Lean Simplificada,  Bookman, 2008 HARRY, M. , LINSENMANND.R., The Six Sigma Fieldbook, Doubleday, New York, 2006KUME, H. (Tradução Miyake, D.I.). 1993. Métodos Estatísticos para Melhoria da Qualidade. São Paulo. Editora Gente, 245 p.HARRY, M. , LINSENMANND.R., The Six Sigma Fieldbook, Doubleday, New York, 2006ISHIKAWA, K. Guide to quality control. Tokyo: Kraus Asian Productivity Organization, 1982.JURAN, J.M.; GRYNA, F.M. Juran controle da qualidade: métodos especiais de apoio à qualidade. São Paulo: Makron Books, 1993LIKER, K. & MEIER D. O Modelo Toyota, Manual de Aplicação. Porto Alegre: Bookman, 2007MARANHÃO, M. ISO Série 9000-Guia de Implementação, Qualitymark, Rio de Janeiro, 2001MONTGOMERY, D.C. 1991. Introduction to Statistical Quality Control. New York. John Wiley & Sons Inc."

# ---------------------------------------------------------------------------
# 5. Column layout: split the merged "A & B" column-width definition (cols
#    1-2 sharing one <col> entry) into a standalone column-1 definition so it
#    matches the target's separate <col> elements for columns 1 and 2.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 30.7109375
